# Update numeric "views/likes" counters (column F) on the "展览" sheet
# and the "全部类型" sheet, as produced by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1409
$ws1.Range("F6").Value  = 71
$ws1.Range("F7").Value  = 11897
$ws1.Range("F8").Value  = 4440
$ws1.Range("F12").Value = 21
$ws1.Range("F17").Value = 5160
$ws1.Range("F19").Value = 195
$ws1.Range("F20").Value = 535
$ws1.Range("F21").Value = 11387
$ws1.Range("F22").Value = 11376

# Sheet 4: "全部类型" (index 4) - same underlying events, rows offset by
# one extra entry inserted earlier in the sheet (row 14).
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 1409
$ws4.Range("F6").Value  = 71
$ws4.Range("F7").Value  = 11897
$ws4.Range("F8").Value  = 4440
$ws4.Range("F12").Value = 21
$ws4.Range("F18").Value = 5160
$ws4.Range("F20").Value = 195
$ws4.Range("F21").Value = 535
$ws4.Range("F22").Value = 11387
$ws4.Range("F23").Value = 11376
